$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.692.52'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '1.583.73'
$ws.Range("E3").Value = '  -3.12%  '
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("E5").Value = '  -2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.506'
$ws.Range("E6").Value = '  -2.57%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.28'
$ws.Range("E8").Value = '  -4.79%  '
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("E10").Value = '  -3.03%  '
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("D12").Value = '1.810.36'
$ws.Range("E12").Value = '  -3.04%  '
$ws.Range("D13").Value = '1.565.93'
$ws.Range("E13").Value = '  -4.25%  '
$ws.Range("E14").Value = '  -3.69%  '
$ws.Range("E15").Value = '  -5.88%  '
$ws.Range("D16").Value = '27.668.80'
$ws.Range("E16").Value = '  -0.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.28'
$ws.Range("E17").Value = '  -2.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.75'
$ws.Range("E18").Value = '  -4.30%  '
$ws.Range("E19").Value = '  -3.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.32'
$ws.Range("E20").Value = '  -6.25%  '
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("E22").Value = '  -4.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.50'
$ws.Range("E24").Value = '  -4.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.15'
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.78'
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.14'
$ws.Range("E28").Value = '  -2.60%  '
$ws.Range("E29").Value = '  -4.05%  '
$ws.Range("E30").Value = '  -2.34%  '
$ws.Range("E31").Value = '  -3.39%  '
$ws.Range("E32").Value = '  -5.28%  '
$ws.Range("D33").Value = '1.388.69'
$ws.Range("E33").Value = '  -0.94%  '
$ws.Range("E34").Value = '  -5.64%  '
$ws.Range("E35").Value = '  -5.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.966'
$ws.Range("E36").Value = '  -4.93%  '
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("E38").Value = '  -3.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.540'
$ws.Range("E39").Value = '  -3.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.820'
$ws.Range("E40").Value = '  -3.77%  '
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("E42").Value = '  -3.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '63.60'
$ws.Range("E43").Value = '  -3.80%  '
$ws.Range("E44").Value = '  +1.57%  '
$ws.Range("E45").Value = '  -4.49%  '
$ws.Range("E46").Value = '  -4.10%  '
$ws.Range("D47").Value = '1.720.08'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.05'
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("E49").Value = '  -2.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0973'
$ws.Range("E50").Value = '  -4.74%  '
$ws.Range("E51").Value = '  -1.39%  '
